$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.747.83'
$ws.Range("E2").Value = '  +0.47%  '

$ws.Range("D3").Value = '1.859.72'
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("E4").Value = '  -1.42%  '

$ws.Range("D5").Value = '320.65'
$ws.Range("E5").Value = '  -0.46%  '

$ws.Range("E6").Value = '  -1.14%  '

$ws.Range("D7").Value = '0.4371'
$ws.Range("E7").Value = '  -0.60%  '

$ws.Range("D8").Value = '0.3776'
$ws.Range("E8").Value = '  -0.41%  '

$ws.Range("D9").Value = '0.07422'
$ws.Range("E9").Value = '  -0.15%  '

$ws.Range("D10").Value = '0.8841'
$ws.Range("E10").Value = '  +0.52%  '

$ws.Range("D11").Value = '21.57'
$ws.Range("E11").Value = '  -0.47%  '

$ws.Range("D12").Value = '1.863.34'
$ws.Range("E12").Value = '  +0.34%  '

$ws.Range("D13").Value = '6.755'
$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("E14").Value = '  -0.74%  '

$ws.Range("D15").Value = '0.07143'
$ws.Range("E15").Value = '  -0.95%  '

$ws.Range("D16").Value = '88.12'

$ws.Range("E17").Value = '  -1.22%  '

$ws.Range("D18").Value = '0.000009045'
$ws.Range("E18").Value = '  -0.31%  '

$ws.Range("E19").Value = '  -1.14%  '

$ws.Range("D20").Value = '15.48'
$ws.Range("E20").Value = '  +0.09%  '

$ws.Range("D21").Value = '27.730.87'
$ws.Range("E21").Value = '  +0.38%  '

$ws.Range("D22").Value = '5.277'
$ws.Range("E22").Value = '  -0.07%  '

$ws.Range("D23").Value = '11.17'
$ws.Range("E23").Value = '  -2.12%  '

$ws.Range("D24").Value = '2.085.49'
$ws.Range("E24").Value = '  +0.30%  '

$ws.Range("E25").Value = '  +5.73%  '

$ws.Range("D26").Value = '157.13'
$ws.Range("E26").Value = '  -0.67%  '

$ws.Range("D27").Value = '18.70'
$ws.Range("E27").Value = '  -0.43%  '

$ws.Range("D28").Value = '5.424'
$ws.Range("E28").Value = '  +2.53%  '

$ws.Range("D29").Value = '1.990'
$ws.Range("E29").Value = '  +0.43%  '

$ws.Range("D30").Value = '121.25'
$ws.Range("E30").Value = '  +3.10%  '

$ws.Range("D31").Value = '0.09058'
$ws.Range("E31").Value = '  -0.12%  '

$ws.Range("E32").Value = '  +0.96%  '

$ws.Range("D33").Value = '0.7716'
$ws.Range("E33").Value = '  +0.92%  '

$ws.Range("D34").Value = '3.036'
$ws.Range("E34").Value = '  +5.04%  '

$ws.Range("D35").Value = '4.571'
$ws.Range("E35").Value = '  +0.69%  '

$ws.Range("D36").Value = '1.019'
$ws.Range("E36").Value = '  -1.08%  '

$ws.Range("D37").Value = '1.137'
$ws.Range("E37").Value = '  -1.43%  '

$ws.Range("D38").Value = '0.01981'
$ws.Range("E38").Value = '  +0.08%  '

$ws.Range("D39").Value = '0.05314'
$ws.Range("E39").Value = '  -0.19%  '

$ws.Range("D40").Value = '2.875'
$ws.Range("E40").Value = '  +1.59%  '

$ws.Range("D41").Value = '0.5176'
$ws.Range("E41").Value = '  +0.06%  '

$ws.Range("D42").Value = '6.962'
$ws.Range("E42").Value = '  +2.43%  '

$ws.Range("E43").Value = '  -0.22%  '

$ws.Range("D44").Value = '8.724'
$ws.Range("E44").Value = '  +2.17%  '

$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '110.28'
$ws.Range("E45").Value = '  +0.91%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '10.78'
$ws.Range("E46").Value = '  +1.66%  '

$ws.Range("D47").Value = '1.718'
$ws.Range("E47").Value = '  +0.03%  '

$ws.Range("D48").Value = '0.4729'
$ws.Range("E48").Value = '  +1.46%  '

$ws.Range("E49").Value = '  -1.25%  '

$ws.Range("D50").Value = '0.06478'
$ws.Range("E50").Value = '  +1.00%  '

$ws.Range("D51").Value = '1.856'
$ws.Range("E51").Value = '  +0.14%  '
